$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($ws, $cellRef, $value)
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextCell $ws 'D2' '29.376.30'
$ws.Range('E2').Value = '  +0.57%  '
Set-TextCell $ws 'D3' '1.874.27'
$ws.Range('E3').Value = '  +0.66%  '
$ws.Range('E4').Value = '  -0.03%  '
Set-TextCell $ws 'D5' '0.7116'
$ws.Range('E5').Value = '  -0.40%  '
Set-TextCell $ws 'D6' '242.03'
$ws.Range('E6').Value = '  +0.74%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('E8').Value = '  +0.96%  '
Set-TextCell $ws 'D9' '0.07780'
$ws.Range('E9').Value = '  +1.18%  '
Set-TextCell $ws 'D10' '25.01'
$ws.Range('E10').Value = '  -0.06%  '
Set-TextCell $ws 'D11' '0.08468'
$ws.Range('E11').Value = '  +1.72%  '
Set-TextCell $ws 'D12' '1.866.47'
$ws.Range('E12').Value = '  -4.01%  '
Set-TextCell $ws 'D13' '5.243'
$ws.Range('E13').Value = '  +0.60%  '
Set-TextCell $ws 'D14' '0.7126'
$ws.Range('E14').Value = '  -0.57%  '
Set-TextCell $ws 'D15' '91.22'
$ws.Range('E15').Value = '  +0.55%  '
Set-TextCell $ws 'D16' '29.379.45'
$ws.Range('E16').Value = '  +0.17%  '
Set-TextCell $ws 'D17' '0.000008242'
$ws.Range('E17').Value = '  +5.62%  '
Set-TextCell $ws 'D18' '6.039'
$ws.Range('E18').Value = '  +1.59%  '
Set-TextCell $ws 'D19' '240.73'
$ws.Range('E19').Value = '  -1.00%  '
Set-TextCell $ws 'D20' '13.25'
$ws.Range('E20').Value = '  +0.81%  '
Set-TextCell $ws 'D21' '2.123.53'
$ws.Range('E21').Value = '  -3.81%  '
Set-TextCell $ws 'D22' '0.9999'
$ws.Range('E22').Value = '  -0.10%  '
Set-TextCell $ws 'D23' '7.806'
$ws.Range('E23').Value = '  -2.41%  '
$ws.Range('E24').Value = '  -0.06%  '
Set-TextCell $ws 'D25' '0.1606'
$ws.Range('E25').Value = '  -0.30%  '
Set-TextCell $ws 'D27' '9.071'
$ws.Range('E27').Value = '  +1.89%  '
Set-TextCell $ws 'D28' '18.48'
$ws.Range('E28').Value = '  -0.49%  '
$ws.Range('E29').Value = '  +1.02%  '
Set-TextCell $ws 'D30' '4.426'
$ws.Range('E30').Value = '  -0.31%  '
$ws.Range('E31').Value = '  +1.79%  '
Set-TextCell $ws 'D32' '1.280'
$ws.Range('E32').Value = '  -4.71%  '
Set-TextCell $ws 'D33' '0.05303'
$ws.Range('E33').Value = '  +2.31%  '
Set-TextCell $ws 'D34' '1.937'
$ws.Range('E34').Value = '  +0.69%  '
Set-TextCell $ws 'D35' '1.177'
$ws.Range('E35').Value = '  +0.54%  '
Set-TextCell $ws 'D36' '0.7469'
$ws.Range('E36').Value = '  -5.77%  '
Set-TextCell $ws 'D37' '2.697'
$ws.Range('E37').Value = '  +0.48%  '
$ws.Range('E38').Value = '  +0.81%  '
$ws.Range('B39').Value = 'Maker'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextCell $ws 'D39' '1.207.67'
$ws.Range('E39').Value = '  +2.08%  '
$ws.Range('B40').Value = 'MXToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextCell $ws 'D40' '2.721'
$ws.Range('E40').Value = '  +1.03%  '
Set-TextCell $ws 'D41' '6.451'
$ws.Range('E41').Value = '  +3.26%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextCell $ws 'D42' '0.8880'
$ws.Range('E42').Value = '  -1.85%  '
$ws.Range('B43').Value = 'Aave'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextCell $ws 'D43' '72.93'
$ws.Range('E43').Value = '  -0.08%  '
Set-TextCell $ws 'D44' '108.14'
$ws.Range('E44').Value = '  +5.59%  '
Set-TextCell $ws 'D45' '0.9999'
$ws.Range('E45').Value = '  +0.01%  '
Set-TextCell $ws 'D46' '2.022.04'
$ws.Range('E46').Value = '  -2.25%  '
$ws.Range('E47').Value = '  +2.39%  '
Set-TextCell $ws 'D48' '0.5208'
$ws.Range('E48').Value = '  +0.05%  '
Set-TextCell $ws 'D49' '0.00000000122'
$ws.Range('E49').Value = '  +7.26%  '
Set-TextCell $ws 'D50' '9.393'
$ws.Range('E50').Value = '  +0.57%  '
Set-TextCell $ws 'D51' '0.4321'
$ws.Range('E51').Value = '  +0.93%  '
